$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) " (complexity levels)" -> " (complexity)"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(" (complexity levels)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " (complexity)", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Insert a new run " using a compute shader" right after the run
#    that contains "modern graphics hardware".
# ---------------------------------------------------------------------
$rngHw = $d.Content.Duplicate
$rngHw.Find.Execute("modern graphics hardware", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
$rngHw.Collapse(0)
$rngHw.InsertAfter(" using a compute shader")

# ---------------------------------------------------------------------
# 3 & 4) Move the "_GoBack" bookmark: remove it from after "Alex Stiyer"
#    (adding a new "_GoBack" bookmark automatically removes the old one,
#    since it is a single-instance bookmark) and place it, collapsed,
#    right after the Mandelbulb Wikipedia hyperlink -- replacing the
#    trailing single-space run that used to sit there.
# ---------------------------------------------------------------------
$rngLink = $d.Content.Duplicate
$rngLink.Find.Execute("https://en.wikipedia.org/wiki/Mandelbulb", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
$hlEnd = $rngLink.End

$spaceRng = $d.Range($hlEnd + 1, $hlEnd + 2)
$d.Bookmarks.Add("_GoBack", $spaceRng)
$spaceRng.Text = ""
